$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Failed : You typed an invalid time.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-01-50-936Z.png'
$ws.Range("G3").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-02-08-853Z.png'
$ws.Range("G4").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-02-26-068Z.png'
$ws.Range("G5").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-02-44-429Z.png'
$ws.Range("G7").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-03-09-236Z.png'
$ws.Range("G8").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-03-27-450Z.png'
$ws.Range("G9").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-03-45-634Z.png'
$ws.Range("G10").Value = ''
$ws.Range("G11").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-04-39-112Z.png'
$ws.Range("G12").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-04-57-164Z.png'
$ws.Range("G13").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-05-15-951Z.png'
$ws.Range("G14").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-05-34-395Z.png'
$ws.Range("G15").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-05-52-673Z.png'
$ws.Range("G16").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T12-06-14-981Z.png'
$ws.Range("G17").Value = ''
$ws.Range("G18").Value = ''
$ws.Range("G19").Value = ''
$ws.Range("G20").Value = ''
$ws.Range("G21").Value = ''
$ws.Range("G22").Value = ''
$ws.Range("G23").Value = ''
$ws.Range("G24").Value = ''
$ws.Range("G27").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T10-33-42-879Z.png'
$ws.Range("G30").Value = 'Failed :  No valid entries or test data issue.'
$ws.Range("G33").Value = 'Failed : [object Promise]& find failed screenshot --> H:\WFM\WFMFailedScreenShot\2024-11-27T10-34-21-802Z.png'
$ws.Range("G37").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T10-34-58-837Z.png'
$ws.Range("G45").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T10-36-01-843Z.png'
$ws.Range("G47").Value = 'Failed : Error Duplicate Punches are not allowed.& find failed Screenshot Path:->H:\WFM\WFMFailedScreenShot\2024-11-27T10-36-26-841Z.png'
